# Adds a new "2022-Q3" sheet (fund holdings data) positioned between the
# existing "总计" summary sheet and the existing "2021-Q4" sheet, and updates
# the "总计" sheet with a corresponding summary row for "2022-Q3" (pushing
# the existing "2021-Q4" summary row down by one).
#
# This host's Worksheets.Add() assigns a fresh sheetId as
# (max existing sheetId + 1). The target workbook expects the NEW "2022-Q3"
# sheet to receive sheetId=2 while the OLD "2021-Q4" sheet becomes sheetId=3.
# To land on that allocation we delete the old "2021-Q4" sheet first (which
# frees sheetId=2), then re-add both sheets in the desired order so the ids
# fall out as: 总计=1 (untouched), 2022-Q3=2 (new), 2021-Q4=3 (re-created
# with data identical to the sheet it replaces).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$oldQ4 = $wb.Worksheets.Item("2021-Q4")

# A cell that already carries the shared "bold / thin-border / centered"
# formatting (style index 2 in the original file) used for header rows and
# the leading index column. Re-used via Range.Copy so the *exact* style is
# reproduced instead of an Excel-normalised approximation.
$styledSrc = $summary.Range("A2")

# Remember the original "2021-Q4" fund-holdings data so it can be re-created
# identically after the sheet is deleted (deletion is only needed to obtain
# the correct sheetId allocation - see note above).
$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q4Rows = @(
    @(0, "001917", "招商量化精选股票A", "2.43", "93.64", "1.19", "0.0289", 9),
    @(1, "007950", "招商量化精选股票C", "0.47", "93.64", "1.19", "0.0056", 9)
)

$oldQ4.Delete()

# --- Create the new "2022-Q3" sheet right after "总计" ---------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$q3Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q3Rows = @(
    @(0, "009514", "创金合信同顺创业板精选股票C", "0.16", "92.10", "2.13", "0.0034", 9),
    @(1, "009513", "创金合信同顺创业板精选股票A", "0.09", "92.10", "2.13", "0.0019", 9),
    @(2, "005021", "渤海汇金量化汇盈灵活配置混合", "0.01", "83.26", "4.08", "0.0004", 5)
)

function Fill-FundSheet($sheet, $headers, $rows) {
    # Header row (B1:H1).
    for ($i = 0; $i -lt $headers.Count; $i++) {
        $col = $i + 2
        $styledSrc.Copy($sheet.Cells.Item(1, $col))
        $sheet.Cells.Item(1, $col).Value = $headers[$i]
    }

    # Data rows.
    for ($r = 0; $r -lt $rows.Count; $r++) {
        $row = $rows[$r]
        $excelRow = $r + 2

        $styledSrc.Copy($sheet.Cells.Item($excelRow, 1))
        $sheet.Cells.Item($excelRow, 1).Value = $row[0]

        for ($col = 2; $col -le 7; $col++) {
            $cell = $sheet.Cells.Item($excelRow, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $row[$col - 1]
        }

        $sheet.Cells.Item($excelRow, 8).Value = $row[7]
    }
}

Fill-FundSheet $q3 $q3Headers $q3Rows

# --- Re-create the "2021-Q4" sheet right after "2022-Q3" -------------------
$q4 = $wb.Worksheets.Add($null, $q3)
$q4.Name = "2021-Q4"

Fill-FundSheet $q4 $q4Headers $q4Rows

# --- Update the "总计" summary sheet ----------------------------------------
# Existing row 2 (A2=0, B2="2021-Q4", C2=2, D2=0.03) becomes the new
# "2022-Q3" summary row, and a new row 3 is added carrying the data that
# used to live in row 2 (now tagged "2021-Q4", with A3=1).
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.01

$styledSrc.Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.03
